$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Locator Type" column with its "CSS" value
$ws.Range("D1").Value = "Locator Type"
$ws.Range("D2").Value = "CSS"

# Widen column C to fit the new header/content, matching the recorded bestFit width
$ws.Columns.Item(3).ColumnWidth = 25.83

# Move the active selection to D3, as captured in the saved workbook state
$ws.Range("D3").Select()
